$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "SEK" / "DOI" column headers are relabeled to lowercase "sek" / "doi"
# (underlying data for rows 2-8 is untouched).
$ws.Range("C1").Value = "sek"
$ws.Range("D1").Value = "doi"

# Move/leave the active selection on D1, matching the saved view state.
$ws.Range("D1").Select() | Out-Null
